$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$csv = @"
ECs,Cthrc1,Ror2,ECs,2,1,2.4135925,4.827185,0.1006166576043813,0.07006619804477368,1,0.5,0.0552105,0.110421,0.01579120108286444,0.01113815545262342,0.13325564872125,0.533022594885,0.001588857872516507,0.0007804082057969886
ECs,Cthrc1,Ror2,FAPs,2,1,2.4135925,4.827185,0.1006166576043813,0.07006619804477368,3,1,2.655685,7.967055,0.7595739188695422,0.8036360573586609,6.409741398362501,38.458448390175,0.07642578892011485,0.05630772315081304
ECs,Cthrc1,Ror2,Neutro,2,1,2.4135925,4.827185,0.1006166576043813,0.07006619804477368,2,0.6666666666666666,0.26551,0.79653,0.07594065982940451,0.08034590306805892,0.640832944675,3.84499766805,0.00764089536830599,0.005629531956452805
ECs,Cthrc1,Ror2,sCs,2,1,2.4135925,4.827185,0.1006166576043813,0.07006619804477368,2,1,0.519877,1.039754,0.1486942202181889,0.1048798841206565,1.2547712281225,5.01908491249,0.01496111544344399,0.007348534731710835
FAPs,Cthrc1,Ror2,ECs,3,1,19.76703633333333,59.301109,0.8240384930780011,0.8607507786564448,1,0.5,0.0552105,0.110421,0.01579120108286444,0.01113815545262342,1.0913479594815,6.548087756889,0.01301255754421531,0.009587175978642138
FAPs,Cthrc1,Ror2,FAPs,3,1,19.76703633333333,59.301109,0.8240384930780011,0.8607507786564448,3,1,2.655685,7.967055,0.7595739188695422,0.8036360573586609,52.49502188488834,472.455196963995,0.6259181474866095,0.6917303621278628
FAPs,Cthrc1,Ror2,Neutro,3,1,19.76703633333333,59.301109,0.8240384930780011,0.8607507786564448,2,0.6666666666666666,0.26551,0.79653,0.07594065982940451,0.08034590306805892,5.248345816863333,47.23511235176999,0.06257802688917159,0.06915779862768695
FAPs,Cthrc1,Ror2,sCs,3,1,19.76703633333333,59.301109,0.8240384930780011,0.8607507786564448,2,1,0.519877,1.039754,0.1486942202181889,0.1048798841206565,10.27642754786433,61.658565287186,0.1225297611580048,0.09027544192225281
Neutro,Cthrc1,Ror2,ECs,3,1,1.151594333333333,3.454783,0.04800709843776271,0.05014589452853316,1,0.5,0.0552105,0.110421,0.01579120108286444,0.01113815545262342,0.06358009894050001,0.381480593643,0.0007580897448355783,0.0005585327685696606
Neutro,Cthrc1,Ror2,FAPs,3,1,1.151594333333333,3.454783,0.04800709843776271,0.05014589452853316,3,1,2.655685,7.967055,0.7595739188695422,0.8036360573586609,3.058271797118334,27.524446174065,0.0364649398939273,0.04029904897163364
Neutro,Cthrc1,Ror2,Neutro,3,1,1.151594333333333,3.454783,0.04800709843776271,0.05014589452853316,2,0.6666666666666666,0.26551,0.79653,0.07594065982940451,0.08034590306805892,0.3057598114433333,2.75183830299,0.003645690731858874,0.004029017181050631
Neutro,Cthrc1,Ror2,sCs,3,1,1.151594333333333,3.454783,0.04800709843776271,0.05014589452853316,2,1,0.519877,1.039754,0.1486942202181889,0.1048798841206565,0.5986874072303334,3.592124443382,0.00713837806714096,0.005259295607279222
sCs,Cthrc1,Ror2,ECs,2,1,0.655778,1.311556,0.02733775087985481,0.01903712877024833,1,0.5,0.0552105,0.110421,0.01579120108286444,0.01113815545262342,0.036205831269,0.144823325076,0.0004316959212970417,0.0002120384996146357
sCs,Cthrc1,Ror2,FAPs,2,1,0.655778,1.311556,0.02733775087985481,0.01903712877024833,3,1,2.655685,7.967055,0.7595739188695422,0.8036360573586609,1.74153979793,10.44923878758,0.0207650425688906,0.0152989231083515
sCs,Cthrc1,Ror2,Neutro,2,1,0.655778,1.311556,0.02733775087985481,0.01903712877024833,2,0.6666666666666666,0.26551,0.79653,0.07594065982940451,0.08034590306805892,0.17411561678,1.04469370068,0.002076046840068058,0.001529555302868528
sCs,Cthrc1,Ror2,sCs,2,1,0.655778,1.311556,0.02733775087985481,0.01903712877024833,2,1,0.519877,1.039754,0.1486942202181889,0.1048798841206565,0.340923899306,1.363695597224,0.004064965549599119,0.001996611859413661
"@

$lines = $csv -split "`n" | Where-Object { $_.Trim().Length -gt 0 }
$rowCount = $lines.Count
$colCount = 20

$data = New-Object 'object[,]' $rowCount,$colCount

for ($r = 0; $r -lt $rowCount; $r++) {
    $fields = $lines[$r].Split(",")
    for ($c = 0; $c -lt $colCount; $c++) {
        $field = $fields[$c]
        if ($c -lt 4) {
            $data[$r, $c] = $field
        } else {
            $data[$r, $c] = [double]$field
        }
    }
}

$startRow = 2
$endRow = $startRow + $rowCount - 1
$endColLetter = [char](64 + $colCount)
$rangeAddr = "A" + $startRow + ":" + $endColLetter + $endRow

$range = $ws.Range($rangeAddr)
$range.Value = $data
